# Update the date heading paragraph
$d = $word.ActiveDocument

$dateFound = $d.Content.Find.Execute("2025-12-05 Friday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-12-06 Saturday", 2)

# Update each math-problem cell in the table, by row/column position
# (cell text uses Range.Text so paragraph/run formatting is preserved).
$t = $d.Tables.Item(1)

$values = @(
    "20-20=",
    "69-11=",
    "84-2=",
    "27+49=",
    "78-57=",
    "84-36=",
    "61+19=",
    "2+69=",
    "11-1=",
    "68-7=",
    "98-45=",
    "85-31=",
    "44+18=",
    "17+24=",
    "76-69=",
    "71-58=",
    "90-8=",
    "81-36=",
    "8+14=",
    "11+11=",
    "34+60=",
    "34-13=",
    "48-29=",
    "56+30=",
    "99-39=",
    "69-33=",
    "17-16=",
    "36+57=",
    "53-13=",
    "61-57=",
    "51-20=",
    "88-66=",
    "6+61=",
    "83-57=",
    "91-45=",
    "21+21=",
    "48-28=",
    "27+22=",
    "99-40=",
    "57-53=",
    "49-25=",
    "76-54=",
    "2+54=",
    "47-28=",
    "21+76=",
    "81-48=",
    "67+6=",
    "59-19=",
    "21+9=",
    "97-84=",
    "59-45=",
    "88+1=",
    "98-74=",
    "45+2=",
    "54+1=",
    "50-20=",
    "78+11=",
    "51-29=",
    "93-25=",
    "97-73=",
    "47-25=",
    "28+40=",
    "27+72=",
    "72-60=",
    "88-41=",
    "34-11=",
    "76-43=",
    "36+16=",
    "62-40=",
    "60-18=",
    "17+76=",
    "60-1=",
    "90-13=",
    "75-30=",
    "77-76=",
    "25+66=",
    "70-32=",
    "22+68=",
    "30+18=",
    "37+17=",
    "66-22=",
    "32+67=",
    "68-23=",
    "78-9=",
    "40+52=",
    "96-69=",
    "83-5=",
    "72-68=",
    "0+68=",
    "22+59=",
    "6+88=",
    "26+60=",
    "63-20=",
    "19+48=",
    "21+4=",
    "98-41=",
    "16+53=",
    "61-28=",
    "32-25=",
    "44+5="
)

$rows = 20
$cols = 5
$idx = 0
for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $cell = $t.Cell($r, $c)
        $cell.Range.Text = $values[$idx]
        $idx++
    }
}

Write-Host "Date replaced:" $dateFound
Write-Host "Cells updated:" $idx
